$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting existing data down one row.
$ws.Rows.Item(1).Insert()

# Add header labels in the new first row (no special style / alignment).
$ws.Range("C1").Value = "주소"
$ws.Range("D1").Value = "이름"

# Re-apply the formatting to F1 that the row-insert left unstyled (match F2's style).
$ws.Range("F2").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Update the active selection to match the saved view state.
$ws.Range("F2").Select()
